{"js": "// The document has a single table; five of its rows (0, 4, 8, 12, 16)\n// hold the \"two-digit \u00f7 one-digit\" practice problems (the rows between\n// them are blank spacer rows). Each diff hunk replaces one cell's text\n// in place, so address every cell positionally via getCell(row, col)\n// rather than matching on old text (several old/new strings repeat).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\ntable.getCell(0, 0).value = \"43\u00f73=14, 1\";\ntable.getCell(0, 1).value = \"98\u00f73=32, 2\";\ntable.getCell(0, 2).value = \"97\u00f72=48, 1\";\ntable.getCell(0, 3).value = \"25\u00f73=8, 1\";\ntable.getCell(0, 4).value = \"87\u00f78=10, 7\";\ntable.getCell(4, 0).value = \"51\u00f73=17, 0\";\ntable.getCell(4, 1).value = \"69\u00f75=13, 4\";\ntable.getCell(4, 2).value = \"41\u00f75=8, 1\";\ntable.getCell(4, 3).value = \"31\u00f72=15, 1\";\ntable.getCell(4, 4).value = \"10\u00f75=2, 0\";\ntable.getCell(8, 0).value = \"61\u00f75=12, 1\";\ntable.getCell(8, 1).value = \"87\u00f76=14, 3\";\ntable.getCell(8, 2).value = \"33\u00f78=4, 1\";\ntable.getCell(8, 3).value = \"28\u00f78=3, 4\";\ntable.getCell(8, 4).value = \"10\u00f75=2, 0\";\ntable.getCell(12, 0).value = \"89\u00f78=11, 1\";\ntable.getCell(12, 1).value = \"38\u00f79=4, 2\";\ntable.getCell(12, 2).value = \"54\u00f77=7, 5\";\ntable.getCell(12, 3).value = \"55\u00f76=9, 1\";\ntable.getCell(12, 4).value = \"25\u00f77=3, 4\";\ntable.getCell(16, 0).value = \"38\u00f77=5, 3\";\ntable.getCell(16, 1).value = \"33\u00f73=11, 0\";\ntable.getCell(16, 2).value = \"11\u00f73=3, 2\";\ntable.getCell(16, 3).value = \"77\u00f74=19, 1\";\ntable.getCell(16, 4).value = \"50\u00f74=12, 2\";\n\nawait context.sync();", "ps1": "# The document has a single table; five of its rows (1, 5, 9, 13, 17 in\n# COM's 1-based indexing) hold the \"two-digit \u00f7 one-digit\" practice\n# problems (the rows between them are blank spacer rows). Each diff hunk\n# replaces one cell's text in place, so address every cell positionally\n# via Cell(row, col) rather than matching on old text (several old/new\n# strings repeat, so a blind global find/replace would be unsafe).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"43\u00f73=14, 1\"\n$t.Cell(1, 2).Range.Text = \"98\u00f73=32, 2\"\n$t.Cell(1, 3).Range.Text = \"97\u00f72=48, 1\"\n$t.Cell(1, 4).Range.Text = \"25\u00f73=8, 1\"\n$t.Cell(1, 5).Range.Text = \"87\u00f78=10, 7\"\n$t.Cell(5, 1).Range.Text = \"51\u00f73=17, 0\"\n$t.Cell(5, 2).Range.Text = \"69\u00f75=13, 4\"\n$t.Cell(5, 3).Range.Text = \"41\u00f75=8, 1\"\n$t.Cell(5, 4).Range.Text = \"31\u00f72=15, 1\"\n$t.Cell(5, 5).Range.Text = \"10\u00f75=2, 0\"\n$t.Cell(9, 1).Range.Text = \"61\u00f75=12, 1\"\n$t.Cell(9, 2).Range.Text = \"87\u00f76=14, 3\"\n$t.Cell(9, 3).Range.Text = \"33\u00f78=4, 1\"\n$t.Cell(9, 4).Range.Text = \"28\u00f78=3, 4\"\n$t.Cell(9, 5).Range.Text = \"10\u00f75=2, 0\"\n$t.Cell(13, 1).Range.Text = \"89\u00f78=11, 1\"\n$t.Cell(13, 2).Range.Text = \"38\u00f79=4, 2\"\n$t.Cell(13, 3).Range.Text = \"54\u00f77=7, 5\"\n$t.Cell(13, 4).Range.Text = \"55\u00f76=9, 1\"\n$t.Cell(13, 5).Range.Text = \"25\u00f77=3, 4\"\n$t.Cell(17, 1).Range.Text = \"38\u00f77=5, 3\"\n$t.Cell(17, 2).Range.Text = \"33\u00f73=11, 0\"\n$t.Cell(17, 3).Range.Text = \"11\u00f73=3, 2\"\n$t.Cell(17, 4).Range.Text = \"77\u00f74=19, 1\"\n$t.Cell(17, 5).Range.Text = \"50\u00f74=12, 2\""}
